# Refresh market price data (columns H-N: currentAveragePrice*, LevePrice*, LeveProfit*)
# across the per-job leve profit worksheets, as produced by the scheduled data-pull runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 350.75
$ws.Range("I33").Value = 311.83334
$ws.Range("J33").Value = 467.5
$ws.Range("K33").Value = 311.83334
$ws.Range("L33").Value = 467.5
$ws.Range("M33").Value = -82.83334000000002
$ws.Range("N33").Value = -925.5

$ws.Range("H86").Value = 2645
$ws.Range("I86").Value = 1931
$ws.Range("J86").Value = 3002
$ws.Range("K86").Value = 1931
$ws.Range("L86").Value = 3002
$ws.Range("M86").Value = -808
$ws.Range("N86").Value = -5248

$ws.Range("H89").Value = 2645
$ws.Range("I89").Value = 1931
$ws.Range("J89").Value = 3002
$ws.Range("K89").Value = 9655
$ws.Range("L89").Value = 15010
$ws.Range("M89").Value = -4039
$ws.Range("N89").Value = -26242

$ws.Range("H121").Value = 880.9375
$ws.Range("I121").Value = 497.5
$ws.Range("J121").Value = 935.7143
$ws.Range("K121").Value = 1492.5
$ws.Range("L121").Value = 2807.1429
$ws.Range("M121").Value = 254.5
$ws.Range("N121").Value = -6301.1429

$ws.Range("H138").Value = 4880316
$ws.Range("I138").Value = 2075.5
$ws.Range("J138").Value = 7409774
$ws.Range("K138").Value = 6226.5
$ws.Range("L138").Value = 22229322
$ws.Range("M138").Value = -1086.5
$ws.Range("N138").Value = -22239602

$ws.Range("H141").Value = 3406.9614
$ws.Range("I141").Value = 3027.6667
$ws.Range("K141").Value = 9083.000100000001
$ws.Range("M141").Value = -3903.000100000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2323.5
$ws.Range("I122").Value = 1784.4
$ws.Range("J122").Value = 2862.6
$ws.Range("K122").Value = 5353.200000000001
$ws.Range("L122").Value = 8587.799999999999
$ws.Range("M122").Value = -2903.200000000001
$ws.Range("N122").Value = -13487.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 11627.167
$ws.Range("I86").Value = 15014.765
$ws.Range("K86").Value = 15014.765
$ws.Range("M86").Value = -13891.765

$ws.Range("H89").Value = 11627.167
$ws.Range("I89").Value = 15014.765
$ws.Range("K89").Value = 75073.825
$ws.Range("M89").Value = -69457.825

$ws.Range("H94").Value = 1006.63635
$ws.Range("I94").Value = 981.1429000000001
$ws.Range("J94").Value = 1051.25
$ws.Range("K94").Value = 981.1429000000001
$ws.Range("L94").Value = 1051.25
$ws.Range("M94").Value = -530.1429000000001
$ws.Range("N94").Value = -1953.25

$ws.Range("H107").Value = 3674.5
$ws.Range("I107").Value = 3010.2856
$ws.Range("J107").Value = 5999.25
$ws.Range("K107").Value = 3010.2856
$ws.Range("L107").Value = 5999.25
$ws.Range("M107").Value = -1090.2856
$ws.Range("N107").Value = -9839.25

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 428.75555
$ws.Range("I107").Value = 416.4074
$ws.Range("J107").Value = 447.27777
$ws.Range("K107").Value = 416.4074
$ws.Range("L107").Value = 447.27777
$ws.Range("M107").Value = 1503.5926
$ws.Range("N107").Value = -4287.27777

$ws.Range("H115").Value = 29000
$ws.Range("J115").Value = 29000
$ws.Range("L115").Value = 29000
$ws.Range("N115").Value = -31350

$ws.Range("H132").Value = 26857.125
$ws.Range("I132").Value = 1605.6666
$ws.Range("K132").Value = 4816.9998
$ws.Range("M132").Value = -2286.9998

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H97").Value = 7235.6665
$ws.Range("I97").Value = 10401.5
$ws.Range("J97").Value = 904
$ws.Range("K97").Value = 31204.5
$ws.Range("L97").Value = 2712
$ws.Range("M97").Value = -30708.5
$ws.Range("N97").Value = -3704

$ws.Range("H131").Value = 919.6774
$ws.Range("I131").Value = 335.7143
$ws.Range("J131").Value = 1090
$ws.Range("K131").Value = 1007.1429
$ws.Range("L131").Value = 3270
$ws.Range("M131").Value = 4032.8571
$ws.Range("N131").Value = -13350

$ws.Range("H136").Value = 3227.1765
$ws.Range("I136").Value = 2475.5715
$ws.Range("J136").Value = 3753.3
$ws.Range("K136").Value = 7426.7145
$ws.Range("L136").Value = 11259.9
$ws.Range("M136").Value = -2326.7145
$ws.Range("N136").Value = -21459.9

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 0
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("M52").ClearContents()
$ws.Range("N52").ClearContents()

$ws.Range("H80").Value = 4647.1665
$ws.Range("J80").Value = 4666.6
$ws.Range("L80").Value = 4666.6
$ws.Range("N80").Value = -6662.6

$ws.Range("H83").Value = 4647.1665
$ws.Range("J83").Value = 4666.6
$ws.Range("L83").Value = 23333
$ws.Range("N83").Value = -33317

$ws.Range("H102").Value = 1724.3889
$ws.Range("I102").Value = 1753
$ws.Range("J102").Value = 1688.625
$ws.Range("K102").Value = 1753
$ws.Range("L102").Value = 1688.625
$ws.Range("M102").Value = -131
$ws.Range("N102").Value = -4932.625

$ws.Range("H107").Value = 236.27777
$ws.Range("J107").Value = 0
$ws.Range("L107").Value = 0
$ws.Range("N107").ClearContents()

$ws.Range("H122").Value = 2999.0908
$ws.Range("I122").Value = 2441.4285
$ws.Range("K122").Value = 7324.2855
$ws.Range("M122").Value = -4874.2855

$ws.Range("H132").Value = 88789.35000000001
$ws.Range("I132").Value = 54087.21
$ws.Range("J132").Value = 253624.5
$ws.Range("K132").Value = 162261.63
$ws.Range("L132").Value = 760873.5
$ws.Range("M132").Value = -159731.63
$ws.Range("N132").Value = -765933.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1358.4546
$ws.Range("I82").Value = 1066.8572
$ws.Range("J82").Value = 1868.75
$ws.Range("K82").Value = 1066.8572
$ws.Range("L82").Value = 1868.75
$ws.Range("M82").Value = -705.8571999999999
$ws.Range("N82").Value = -2590.75

$ws.Range("H85").Value = 1358.4546
$ws.Range("I85").Value = 1066.8572
$ws.Range("J85").Value = 1868.75
$ws.Range("K85").Value = 1066.8572
$ws.Range("L85").Value = 1868.75
$ws.Range("M85").Value = 181.1428000000001
$ws.Range("N85").Value = -4364.75

$ws.Range("H93").Value = 1088.52
$ws.Range("I93").Value = 1088.52
$ws.Range("K93").Value = 1088.52
$ws.Range("M93").Value = 159.48

$ws.Range("H117").Value = 49690
$ws.Range("J117").Value = 49690
$ws.Range("L117").Value = 49690
$ws.Range("N117").Value = -58868

$ws.Range("H122").Value = 3733.3235
$ws.Range("I122").Value = 4363.857
$ws.Range("J122").Value = 3291.95
$ws.Range("K122").Value = 13091.571
$ws.Range("L122").Value = 9875.849999999999
$ws.Range("M122").Value = -10641.571
$ws.Range("N122").Value = -14775.85

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 451.6216
$ws.Range("I107").Value = 366.5
$ws.Range("J107").Value = 652.8182
$ws.Range("K107").Value = 1099.5
$ws.Range("L107").Value = 1958.4546
$ws.Range("M107").Value = 820.5
$ws.Range("N107").Value = -5798.4546

$ws.Range("H118").Value = 40345.5
$ws.Range("J118").Value = 40345.5
$ws.Range("L118").Value = 40345.5
$ws.Range("N118").Value = -43659.5

Write-Host "Updated cells across 8 worksheets."
